$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add tests for scalar_annualized (Scalar_annualized2 / discrete value, Scalar_annualized3 / log value & std)
# Order of writes chosen to reproduce the author's shared-string insertion order.
$ws.Range("A38").Value = "Scalar_annualized2"
$ws.Range("C38").Value = "scalar_annualized_test2"
$ws.Range("B37").Value = "Test scalar annulized for discrete value"
$ws.Range("A39").Value = "Scalar_annualized3"
$ws.Range("B38").Value = "Test scalar annulized for log value"
$ws.Range("C39").Value = "scalar_annualized_test3"
$ws.Range("B39").Value = "Test scalar annulized for std"

# Match the author's final selection state
$ws.Range("B39").Select()
